$wb = $excel.ActiveWorkbook

# The "Repayment Schedule" worksheet gets a new blank column inserted
# before column N (pushing old N..P to O..Q).
$ws = $wb.Worksheets.Item("Repayment Schedule")

$ws.Columns("N:N").Insert()

# Select the new sheet/range the way the author left it in the file.
$ws.Range("R17").Select()

# The previously active sheet ("Loan Tranche Details") loses the
# tabSelected flag, and "Repayment Schedule" becomes the active tab.
$ws.Activate()

$wb.Save()
